# Apply the 2021-03-22 -> 2021-03-23 "as of" date update to the confidential
# disclosure cell, and refresh the Weight / Percent Change figures in
# columns D and E (rows 2-15) to match the newly re-run model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet ships sheet-protected; temporarily unprotect so the cells
# below can be written, then restore protection once done.
$ws.Unprotect()

# --- Update the "Model holdings provided as of ..." disclosure text -------
$disclosure = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."
$ws.Range("A18").Value = $disclosure

# --- Refresh Weight (D) / Percent Change (E) figures -----------------------
$ws.Range("D2").Value  = 0.05579350697891961
$ws.Range("E2").Value  = -0.01138731404758442

$ws.Range("D3").Value  = 0.02343384068016333
$ws.Range("E3").Value  = -0.01514233797698361

$ws.Range("D4").Value  = 0.03201496153838355
$ws.Range("E4").Value  = -0.0163404902147064

$ws.Range("D5").Value  = 0.03178380658142771
$ws.Range("E5").Value  = -0.02222222222222214

$ws.Range("D6").Value  = 0.0333917529048135
$ws.Range("E6").Value  = -0.03473402963322814

$ws.Range("D7").Value  = 0.01892023599434144
$ws.Range("E7").Value  = -0.02788697174293575

$ws.Range("D8").Value  = 0.00484614339582853
$ws.Range("E8").Value  = -0.0188284518828451

$ws.Range("D9").Value  = 0.006522354779601859
$ws.Range("E9").Value  = -0.007875647668393881

$ws.Range("D10").Value = 0.07007645010871719
$ws.Range("E10").Value = -0.006944444444444753

$ws.Range("D11").Value = 0.07023866411359847
$ws.Range("E11").Value = -0.006928406466512715

$ws.Range("D12").Value = 0.1473552020341636
$ws.Range("E12").Value = 0.009026860413914761

$ws.Range("D13").Value = 0.390895198262688
$ws.Range("E13").Value = 0.004938271604938205

$ws.Range("D14").Value = 0.1147278826273532
$ws.Range("E14").Value = 0.009526166027465033

# Row 15 (Total) keeps its Weight of 1 (100%); only Percent Change changes.
$ws.Range("E15").Value = -0.0006695653408150815

# Restore sheet protection.
$ws.Protect()
